# Edit script: add "aws.ses" command group (sendMail / sendTextMail) to the
# hidden '#system' sheet, plus two new commands (base64 under "io", upload
# under "ws"), and register a new "aws.ses" defined name.  This replays the
# manual spreadsheet edit described by the commit message / xml diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Insert a brand-new column C ("aws.ses") - this pushes every existing
#    column from C onward (base, csv, desktop, ... xml) one slot to the
#    right (C->D, D->E, ... Z->AA).
# ---------------------------------------------------------------------
$ws.Columns("C").Insert()

$ws.Cells.Item(1, 3).Value = "aws.ses"
$ws.Cells.Item(2, 3).Value = "sendMail(profile,to,subject,body)"
$ws.Cells.Item(3, 3).Value = "sendTextMail(profile,to,subject,body)"

# ---------------------------------------------------------------------
# 2) Column A ("target") lists the name of every command group; insert a
#    new entry "aws.ses" right after "aws.s3" (row 3), shifting every row
#    below it down by one. We shift manually (bottom-up) instead of using
#    Range.Insert because Insert() on this runtime shifts the *entire*
#    row, not just the target column.
# ---------------------------------------------------------------------
for ($r = 26; $r -ge 3; $r--) {
    $srcVal = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r + 1, 1).Value = $srcVal
}
$ws.Cells.Item(3, 1).Value = "aws.ses"

# ---------------------------------------------------------------------
# 3) Column J ("io", formerly column I before the column insert above)
#    gains a new command "base64(var,file)" inserted alphabetically
#    between "assertReadableFile(file,minByte)" (row 4) and
#    "compare(expected,actual,failFast)" (row 5); shift rows 5-23 down to
#    6-24.
# ---------------------------------------------------------------------
for ($r = 23; $r -ge 5; $r--) {
    $srcVal = $ws.Cells.Item($r, 10).Value2
    $ws.Cells.Item($r + 1, 10).Value = $srcVal
}
$ws.Cells.Item(5, 10).Value = "base64(var,file)"

# ---------------------------------------------------------------------
# 4) Column Y ("ws", formerly column X) gains a new command
#    "upload(url,body,fileParams,var)" appended as the new last row (17).
# ---------------------------------------------------------------------
$ws.Cells.Item(17, 25).Value = "upload(url,body,fileParams,var)"

# ---------------------------------------------------------------------
# 5) Update the defined names so they keep pointing at the right columns
#    now that everything from C onward moved one column to the right, and
#    register the brand-new "aws.ses" name. Ranges whose row-count grew
#    (io, target, ws) get their final row bumped too.
# ---------------------------------------------------------------------
$wb.Names.Item("base").RefersTo        = "='#system'!`$D`$2:`$D`$36"
$wb.Names.Item("csv").RefersTo         = "='#system'!`$E`$2:`$E`$5"
$wb.Names.Item("desktop").RefersTo     = "='#system'!`$F`$2:`$F`$92"
$wb.Names.Item("excel").RefersTo       = "='#system'!`$G`$2:`$G`$14"
$wb.Names.Item("external").RefersTo    = "='#system'!`$H`$2:`$H`$3"
$wb.Names.Item("image").RefersTo       = "='#system'!`$I`$2:`$I`$5"
$wb.Names.Item("io").RefersTo          = "='#system'!`$J`$2:`$J`$24"
$wb.Names.Item("jms").RefersTo         = "='#system'!`$K`$2:`$K`$4"
$wb.Names.Item("json").RefersTo        = "='#system'!`$L`$2:`$L`$14"
$wb.Names.Item("mail").RefersTo        = "='#system'!`$M`$2:`$M`$2"
$wb.Names.Item("number").RefersTo      = "='#system'!`$N`$2:`$N`$15"
$wb.Names.Item("pdf").RefersTo         = "='#system'!`$O`$2:`$O`$16"
$wb.Names.Item("rdbms").RefersTo       = "='#system'!`$P`$2:`$P`$7"
$wb.Names.Item("redis").RefersTo       = "='#system'!`$Q`$2:`$Q`$10"
$wb.Names.Item("ssh").RefersTo         = "='#system'!`$T`$2:`$T`$9"
$wb.Names.Item("step").RefersTo        = "='#system'!`$U`$2:`$U`$4"
$wb.Names.Item("target").RefersTo      = "='#system'!`$A`$2:`$A`$27"
$wb.Names.Item("web").RefersTo         = "='#system'!`$V`$2:`$V`$117"
$wb.Names.Item("webalert").RefersTo    = "='#system'!`$W`$2:`$W`$8"
$wb.Names.Item("webcookie").RefersTo   = "='#system'!`$X`$2:`$X`$8"
$wb.Names.Item("ws").RefersTo          = "='#system'!`$Y`$2:`$Y`$17"
$wb.Names.Item("xml").RefersTo         = "='#system'!`$AA`$2:`$AA`$11"
$wb.Names.Item("sms").RefersTo         = "='#system'!`$R`$2:`$R`$2"
$wb.Names.Item("sound").RefersTo       = "='#system'!`$S`$2:`$S`$5"
$wb.Names.Item("ws.async").RefersTo    = "='#system'!`$Z`$2:`$Z`$8"

$wb.Names.Add("aws.ses", "='#system'!`$C`$2:`$C`$3")

Write-Output "done"
